# Allow for text formatting
#
# The "Action" column's GOTO-style cells in rows 3-5 are rewritten to all
# point at the (new) "GOTO(PreviousMedBrand)" action, and the now-unused
# "LOAD(SomeOtheProcess); JUMP(Dropout)" action is removed entirely.
#
# E4 and E5 are retyped first (plain text) so they end up sharing one
# shared-string entry, and E3 is retyped last with part of its text
# explicitly (re)formatted -- that's what turns it into its own rich-text
# shared-string entry made up of two runs ("GOTO(" + "PreviousMedBrand)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E4 (was "JUMP(SafetyQuestions)") -------------------------------------
$ws.Range("E4").Value = "GOTO(PreviousMedBrand)"
$ws.Range("E4").Font.Name = "Calibri"
$ws.Range("E4").Font.Size = 11
$ws.Range("E4").Font.Color = 0

# --- E5 (was "LOAD(SomeOtheProcess); JUMP(Dropout)") -----------------------
$ws.Range("E5").Value = "GOTO(PreviousMedBrand)"
$ws.Range("E5").Font.Name = "Calibri"
$ws.Range("E5").Font.Size = 11
$ws.Range("E5").Font.Color = 0

# --- E3 (was "IF [previous_var] == some_value THEN GOTO(PrevMedication) ---
# --- ELSE GOTO(OtherQ)") ----------------------------------------------------
$ws.Range("E3").Value = "GOTO(PreviousMedBrand)"

# Re-apply explicit character formatting to only the new jump target, which
# is what splits the text into the two runs that show up in the saved
# file: the literal "GOTO(" (left on the sheet's base font) and
# "PreviousMedBrand)" (explicitly reformatted).
$suffix = $ws.Range("E3").Characters(6, 17)
$suffix.Font.Name = "Calibri"
$suffix.Font.Size = 11
$suffix.Font.Color = 0

# --- Row heights reflow slightly once the new formatting is in place -------
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(4).RowHeight = 14.9
$ws.Rows.Item(5).RowHeight = 14.9

# --- Leave the selection / scroll position where editing ended up ----------
$ws.Range("C17").Select()
$excel.ActiveWindow.ScrollColumn = 3
